$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Step 1: Insert a new "2022-Q3" worksheet right after "总计",
# by duplicating the structurally-identical "2022-Q2" sheet so
# formatting/styles match, then overwrite its values.
# ---------------------------------------------------------------
$summarySheet = $wb.Worksheets.Item("总计")
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($null, $summarySheet)
$new = $wb.Worksheets.Item(2)
$new.Name = "2022-Q3"

# Extend column A styling (bold/border) down to the two extra rows (36-37)
$new.Cells.Item(35,1).Copy()
$new.Range("A36:A37").PasteSpecial(-4122)

$new.Cells.Item(2,1).Value = 0
$new.Cells.Item(2,2).Value = '''510810'
$new.Cells.Item(2,3).Value = '汇添富中证上海国企ETF'
$new.Cells.Item(2,4).Value = '''63.53'
$new.Cells.Item(2,5).Value = '''98.46'
$new.Cells.Item(2,6).Value = '''7.13'
$new.Cells.Item(2,7).Value = '''4.5297'
$new.Cells.Item(2,8).Value = 4
$new.Cells.Item(3,1).Value = 1
$new.Cells.Item(3,2).Value = '''004854'
$new.Cells.Item(3,3).Value = '广发中证全指汽车指数A'
$new.Cells.Item(3,4).Value = '''15.34'
$new.Cells.Item(3,5).Value = '''94.07'
$new.Cells.Item(3,6).Value = '''13.26'
$new.Cells.Item(3,7).Value = '''2.0341'
$new.Cells.Item(3,8).Value = 2
$new.Cells.Item(4,1).Value = 2
$new.Cells.Item(4,2).Value = '''004855'
$new.Cells.Item(4,3).Value = '广发中证全指汽车指数C'
$new.Cells.Item(4,4).Value = '''10.15'
$new.Cells.Item(4,5).Value = '''94.07'
$new.Cells.Item(4,6).Value = '''13.26'
$new.Cells.Item(4,7).Value = '''1.3459'
$new.Cells.Item(4,8).Value = 2
$new.Cells.Item(5,1).Value = 3
$new.Cells.Item(5,2).Value = '''516110'
$new.Cells.Item(5,3).Value = '国泰中证800汽车与零部件ETF'
$new.Cells.Item(5,4).Value = '''5.83'
$new.Cells.Item(5,5).Value = '''97.46'
$new.Cells.Item(5,6).Value = '''11.33'
$new.Cells.Item(5,7).Value = '''0.6605'
$new.Cells.Item(5,8).Value = 2
$new.Cells.Item(6,1).Value = 4
$new.Cells.Item(6,2).Value = '''510150'
$new.Cells.Item(6,3).Value = '招商上证消费80ETF'
$new.Cells.Item(6,4).Value = '''11.46'
$new.Cells.Item(6,5).Value = '''99.70'
$new.Cells.Item(6,6).Value = '''2.48'
$new.Cells.Item(6,7).Value = '''0.2842'
$new.Cells.Item(6,8).Value = 10
$new.Cells.Item(7,1).Value = 5
$new.Cells.Item(7,2).Value = '''090007'
$new.Cells.Item(7,3).Value = '大成策略回报混合'
$new.Cells.Item(7,4).Value = '''9.90'
$new.Cells.Item(7,5).Value = '''61.86'
$new.Cells.Item(7,6).Value = '''2.42'
$new.Cells.Item(7,7).Value = '''0.2396'
$new.Cells.Item(7,8).Value = 10
$new.Cells.Item(8,1).Value = 6
$new.Cells.Item(8,2).Value = '''012069'
$new.Cells.Item(8,3).Value = '天弘安康颐享12个月持有期混合A'
$new.Cells.Item(8,4).Value = '''20.65'
$new.Cells.Item(8,5).Value = '''20.51'
$new.Cells.Item(8,6).Value = '''1.08'
$new.Cells.Item(8,7).Value = '''0.2230'
$new.Cells.Item(8,8).Value = 2
$new.Cells.Item(9,1).Value = 7
$new.Cells.Item(9,2).Value = '''560003'
$new.Cells.Item(9,3).Value = '益民创新优势混合'
$new.Cells.Item(9,4).Value = '''4.77'
$new.Cells.Item(9,5).Value = '''89.16'
$new.Cells.Item(9,6).Value = '''3.30'
$new.Cells.Item(9,7).Value = '''0.1574'
$new.Cells.Item(9,8).Value = 8
$new.Cells.Item(10,1).Value = 8
$new.Cells.Item(10,2).Value = '''510160'
$new.Cells.Item(10,3).Value = '南方中证南方小康产业ETF'
$new.Cells.Item(10,4).Value = '''2.34'
$new.Cells.Item(10,5).Value = '''99.43'
$new.Cells.Item(10,6).Value = '''6.37'
$new.Cells.Item(10,7).Value = '''0.1491'
$new.Cells.Item(10,8).Value = 2
$new.Cells.Item(11,1).Value = 9
$new.Cells.Item(11,2).Value = '''011977'
$new.Cells.Item(11,3).Value = '格林研究优选混合A'
$new.Cells.Item(11,4).Value = '''1.96'
$new.Cells.Item(11,5).Value = '''92.69'
$new.Cells.Item(11,6).Value = '''5.44'
$new.Cells.Item(11,7).Value = '''0.1066'
$new.Cells.Item(11,8).Value = 8
$new.Cells.Item(12,1).Value = 10
$new.Cells.Item(12,2).Value = '''013463'
$new.Cells.Item(12,3).Value = '大成致远优势一年持有期混合A'
$new.Cells.Item(12,4).Value = '''3.65'
$new.Cells.Item(12,5).Value = '''60.88'
$new.Cells.Item(12,6).Value = '''2.51'
$new.Cells.Item(12,7).Value = '''0.0916'
$new.Cells.Item(12,8).Value = 10
$new.Cells.Item(13,1).Value = 11
$new.Cells.Item(13,2).Value = '''000646'
$new.Cells.Item(13,3).Value = '华润元大量化优选混合A'
$new.Cells.Item(13,4).Value = '''1.47'
$new.Cells.Item(13,5).Value = '''73.62'
$new.Cells.Item(13,6).Value = '''4.93'
$new.Cells.Item(13,7).Value = '''0.0725'
$new.Cells.Item(13,8).Value = 4
$new.Cells.Item(14,1).Value = 12
$new.Cells.Item(14,2).Value = '''014179'
$new.Cells.Item(14,3).Value = '中银证券远见价值混合A'
$new.Cells.Item(14,4).Value = '''1.56'
$new.Cells.Item(14,5).Value = '''93.65'
$new.Cells.Item(14,6).Value = '''3.59'
$new.Cells.Item(14,7).Value = '''0.0560'
$new.Cells.Item(14,8).Value = 8
$new.Cells.Item(15,1).Value = 13
$new.Cells.Item(15,2).Value = '''159936'
$new.Cells.Item(15,3).Value = '广发中证全指可选消费ETF指数'
$new.Cells.Item(15,4).Value = '''1.82'
$new.Cells.Item(15,5).Value = '''98.30'
$new.Cells.Item(15,6).Value = '''2.34'
$new.Cells.Item(15,7).Value = '''0.0426'
$new.Cells.Item(15,8).Value = 6
$new.Cells.Item(16,1).Value = 14
$new.Cells.Item(16,2).Value = '''010764'
$new.Cells.Item(16,3).Value = '九泰锐升混合'
$new.Cells.Item(16,4).Value = '''2.15'
$new.Cells.Item(16,5).Value = '''54.42'
$new.Cells.Item(16,6).Value = '''1.85'
$new.Cells.Item(16,7).Value = '''0.0398'
$new.Cells.Item(16,8).Value = 10
$new.Cells.Item(17,1).Value = 15
$new.Cells.Item(17,2).Value = '''000042'
$new.Cells.Item(17,3).Value = '中证财通中国可持续发展100 (ECPI ESG) 指数增强A'
$new.Cells.Item(17,4).Value = '''1.73'
$new.Cells.Item(17,5).Value = '''93.55'
$new.Cells.Item(17,6).Value = '''1.89'
$new.Cells.Item(17,7).Value = '''0.0327'
$new.Cells.Item(17,8).Value = 3
$new.Cells.Item(18,1).Value = 16
$new.Cells.Item(18,2).Value = '''005571'
$new.Cells.Item(18,3).Value = '中银证券新能源灵活配置混合A'
$new.Cells.Item(18,4).Value = '''0.53'
$new.Cells.Item(18,5).Value = '''90.32'
$new.Cells.Item(18,6).Value = '''5.65'
$new.Cells.Item(18,7).Value = '''0.0299'
$new.Cells.Item(18,8).Value = 4
$new.Cells.Item(19,1).Value = 17
$new.Cells.Item(19,2).Value = '''159872'
$new.Cells.Item(19,3).Value = '鹏华中证车联网主题ETF'
$new.Cells.Item(19,4).Value = '''0.52'
$new.Cells.Item(19,5).Value = '''97.61'
$new.Cells.Item(19,6).Value = '''4.56'
$new.Cells.Item(19,7).Value = '''0.0237'
$new.Cells.Item(19,8).Value = 8
$new.Cells.Item(20,1).Value = 18
$new.Cells.Item(20,2).Value = '''003980'
$new.Cells.Item(20,3).Value = '中银证券瑞益灵活配置混合A'
$new.Cells.Item(20,4).Value = '''0.66'
$new.Cells.Item(20,5).Value = '''91.21'
$new.Cells.Item(20,6).Value = '''3.39'
$new.Cells.Item(20,7).Value = '''0.0224'
$new.Cells.Item(20,8).Value = 8
$new.Cells.Item(21,1).Value = 19
$new.Cells.Item(21,2).Value = '''011243'
$new.Cells.Item(21,3).Value = '万家惠裕回报6个月持有期混合A'
$new.Cells.Item(21,4).Value = '''1.54'
$new.Cells.Item(21,5).Value = '''27.67'
$new.Cells.Item(21,6).Value = '''1.38'
$new.Cells.Item(21,7).Value = '''0.0213'
$new.Cells.Item(21,8).Value = 2
$new.Cells.Item(22,1).Value = 20
$new.Cells.Item(22,2).Value = '''007807'
$new.Cells.Item(22,3).Value = '建信MSCI中国A股指数增强C'
$new.Cells.Item(22,4).Value = '''1.21'
$new.Cells.Item(22,5).Value = '''92.81'
$new.Cells.Item(22,6).Value = '''1.67'
$new.Cells.Item(22,7).Value = '''0.0202'
$new.Cells.Item(22,8).Value = 10
$new.Cells.Item(23,1).Value = 21
$new.Cells.Item(23,2).Value = '''005572'
$new.Cells.Item(23,3).Value = '中银证券新能源灵活配置混合C'
$new.Cells.Item(23,4).Value = '''0.25'
$new.Cells.Item(23,5).Value = '''90.32'
$new.Cells.Item(23,6).Value = '''5.65'
$new.Cells.Item(23,7).Value = '''0.0141'
$new.Cells.Item(23,8).Value = 4
$new.Cells.Item(24,1).Value = 22
$new.Cells.Item(24,2).Value = '''515500'
$new.Cells.Item(24,3).Value = '海富通中证长三角领先ETF'
$new.Cells.Item(24,4).Value = '''0.30'
$new.Cells.Item(24,5).Value = '''97.12'
$new.Cells.Item(24,6).Value = '''4.53'
$new.Cells.Item(24,7).Value = '''0.0136'
$new.Cells.Item(24,8).Value = 5
$new.Cells.Item(25,1).Value = 23
$new.Cells.Item(25,2).Value = '''007806'
$new.Cells.Item(25,3).Value = '建信MSCI中国A股指数增强A'
$new.Cells.Item(25,4).Value = '''0.80'
$new.Cells.Item(25,5).Value = '''92.81'
$new.Cells.Item(25,6).Value = '''1.67'
$new.Cells.Item(25,7).Value = '''0.0134'
$new.Cells.Item(25,8).Value = 10
$new.Cells.Item(26,1).Value = 24
$new.Cells.Item(26,2).Value = '''012070'
$new.Cells.Item(26,3).Value = '天弘安康颐享12个月持有期混合C'
$new.Cells.Item(26,4).Value = '''1.19'
$new.Cells.Item(26,5).Value = '''20.51'
$new.Cells.Item(26,6).Value = '''1.08'
$new.Cells.Item(26,7).Value = '''0.0129'
$new.Cells.Item(26,8).Value = 2
$new.Cells.Item(27,1).Value = 25
$new.Cells.Item(27,2).Value = '''005083'
$new.Cells.Item(27,3).Value = '诺德量化蓝筹增强混合C'
$new.Cells.Item(27,4).Value = '''0.57'
$new.Cells.Item(27,5).Value = '''92.85'
$new.Cells.Item(27,6).Value = '''2.19'
$new.Cells.Item(27,7).Value = '''0.0125'
$new.Cells.Item(27,8).Value = 10
$new.Cells.Item(28,1).Value = 26
$new.Cells.Item(28,2).Value = '''011978'
$new.Cells.Item(28,3).Value = '格林研究优选混合C'
$new.Cells.Item(28,4).Value = '''0.23'
$new.Cells.Item(28,5).Value = '''92.69'
$new.Cells.Item(28,6).Value = '''5.44'
$new.Cells.Item(28,7).Value = '''0.0125'
$new.Cells.Item(28,8).Value = 8
$new.Cells.Item(29,1).Value = 27
$new.Cells.Item(29,2).Value = '''007827'
$new.Cells.Item(29,3).Value = '华润元大量化优选混合C'
$new.Cells.Item(29,4).Value = '''0.19'
$new.Cells.Item(29,5).Value = '''73.62'
$new.Cells.Item(29,6).Value = '''4.93'
$new.Cells.Item(29,7).Value = '''0.0094'
$new.Cells.Item(29,8).Value = 4
$new.Cells.Item(30,1).Value = 28
$new.Cells.Item(30,2).Value = '''003981'
$new.Cells.Item(30,3).Value = '中银证券瑞益灵活配置混合C'
$new.Cells.Item(30,4).Value = '''0.19'
$new.Cells.Item(30,5).Value = '''91.21'
$new.Cells.Item(30,6).Value = '''3.39'
$new.Cells.Item(30,7).Value = '''0.0064'
$new.Cells.Item(30,8).Value = 8
$new.Cells.Item(31,1).Value = 29
$new.Cells.Item(31,2).Value = '''014180'
$new.Cells.Item(31,3).Value = '中银证券远见价值混合C'
$new.Cells.Item(31,4).Value = '''0.16'
$new.Cells.Item(31,5).Value = '''93.65'
$new.Cells.Item(31,6).Value = '''3.59'
$new.Cells.Item(31,7).Value = '''0.0057'
$new.Cells.Item(31,8).Value = 8
$new.Cells.Item(32,1).Value = 30
$new.Cells.Item(32,2).Value = '''013464'
$new.Cells.Item(32,3).Value = '大成致远优势一年持有期混合C'
$new.Cells.Item(32,4).Value = '''0.17'
$new.Cells.Item(32,5).Value = '''60.88'
$new.Cells.Item(32,6).Value = '''2.51'
$new.Cells.Item(32,7).Value = '''0.0043'
$new.Cells.Item(32,8).Value = 10
$new.Cells.Item(33,1).Value = 31
$new.Cells.Item(33,2).Value = '''011244'
$new.Cells.Item(33,3).Value = '万家惠裕回报6个月持有期混合C'
$new.Cells.Item(33,4).Value = '''0.12'
$new.Cells.Item(33,5).Value = '''27.67'
$new.Cells.Item(33,6).Value = '''1.38'
$new.Cells.Item(33,7).Value = '''0.0017'
$new.Cells.Item(33,8).Value = 2
$new.Cells.Item(34,1).Value = 32
$new.Cells.Item(34,2).Value = '''004695'
$new.Cells.Item(34,3).Value = '东兴未来价值灵活配置混合A'
$new.Cells.Item(34,4).Value = '''0.05'
$new.Cells.Item(34,5).Value = '''82.42'
$new.Cells.Item(34,6).Value = '''3.08'
$new.Cells.Item(34,7).Value = '''0.0015'
$new.Cells.Item(34,8).Value = 10
$new.Cells.Item(35,1).Value = 33
$new.Cells.Item(35,2).Value = '''003184'
$new.Cells.Item(35,3).Value = '中证财通中国可持续发展100 (ECPI ESG) 指数增强C'
$new.Cells.Item(35,4).Value = '''0.00'
$new.Cells.Item(35,5).Value = '''93.55'
$new.Cells.Item(35,6).Value = '''1.89'
$new.Cells.Item(35,7).Value = 0
$new.Cells.Item(35,8).Value = 3
$new.Cells.Item(36,1).Value = 34
$new.Cells.Item(36,2).Value = '''007550'
$new.Cells.Item(36,3).Value = '东兴未来价值灵活配置混合C'
$new.Cells.Item(36,4).Value = '''0.00'
$new.Cells.Item(36,5).Value = '''82.42'
$new.Cells.Item(36,6).Value = '''3.08'
$new.Cells.Item(36,7).Value = 0
$new.Cells.Item(36,8).Value = 10
$new.Cells.Item(37,1).Value = 35
$new.Cells.Item(37,2).Value = '''005082'
$new.Cells.Item(37,3).Value = '诺德量化蓝筹增强混合A'
$new.Cells.Item(37,4).Value = '''0.00'
$new.Cells.Item(37,5).Value = '''92.85'
$new.Cells.Item(37,6).Value = '''2.19'
$new.Cells.Item(37,7).Value = 0
$new.Cells.Item(37,8).Value = 10
# ---------------------------------------------------------------
# Step 2: Update "总计" (summary) sheet: insert a new row 2 for
# 2022-Q3 and shift the rest of the rows down by one.
# ---------------------------------------------------------------
$summarySheet.Rows.Item(2).Insert()

# Re-apply the index-column (A) style from the row below onto the
# newly inserted row, then clear the stray format Excel applied to
# B2:D2 so they match the unstyled data cells elsewhere.
$summarySheet.Cells.Item(3,1).Copy()
$summarySheet.Cells.Item(2,1).PasteSpecial(-4122)
$summarySheet.Range("B2:D2").ClearFormats()

$summarySheet.Cells.Item(2,1).Value = 0
$summarySheet.Cells.Item(2,2).Value = "2022-Q3"
$summarySheet.Cells.Item(2,3).Value = 36
$summarySheet.Cells.Item(2,4).Value = 10.29

# Renumber the index column (A) for the rows that shifted down
for ($i = 3; $i -le 9; $i++) {
    $summarySheet.Cells.Item($i,1).Value = $i - 2
}
